$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1 (rows 31-41, "Bodrumi"/basement vertices): update surveyed values.
# Row count stays the same; only numeric values change and the stray Q/R
# helper cells (columns 17/18) get cleared out.
# ---------------------------------------------------------------------------

$ws.Range("B31").Value = 47
$ws.Range("C31").Value = 7510657.7807
$ws.Range("D31").Value = 4693488.7863999996
$ws.Range("I31").Value = 77.703999999999994

$ws.Range("B32").Value = 48
$ws.Range("C32").Value = 7510657.5647999998
$ws.Range("D32").Value = 4693488.6553999996

$ws.Range("B33").Value = 49
$ws.Range("C33").Value = 7510650.6701999996
$ws.Range("D33").Value = 4693484.4735000003

$ws.Range("B34").Value = 56
$ws.Range("C34").Value = 7510662.5954999998
$ws.Range("D34").Value = 4693480.5250000004

$ws.Range("B35").Value = 57
$ws.Range("C35").Value = 7510658.7253
$ws.Range("D35").Value = 4693478.0266000004

$ws.Range("B36").Value = 58
$ws.Range("C36").Value = 7510656.9850000003
$ws.Range("D36").Value = 4693476.9709999999

$ws.Range("B37").Value = 73
$ws.Range("C37").Value = 7510658.6602999996
$ws.Range("D37").Value = 4693478.1381000001

$ws.Range("B38").Value = 74
$ws.Range("C38").Value = 7510655.1160000004
$ws.Range("D38").Value = 4693477.1179999998
$ws.Range("Q38").ClearContents()
$ws.Range("R38").ClearContents()

$ws.Range("B39").Value = 87
$ws.Range("C39").Value = 7510656.4336000001
$ws.Range("D39").Value = 4693477.9172
$ws.Range("I39").Value = 8.6590000000000007

$ws.Range("B40").Value = 45
$ws.Range("C40").Value = 7510650.1025
$ws.Range("D40").Value = 4693485.4210999999
$ws.Range("Q40").ClearContents()

$ws.Range("B41").Value = 46
$ws.Range("C41").Value = 7510657.0235000001
$ws.Range("D41").Value = 4693489.5478999997
$ws.Range("Q41").ClearContents()

# ---------------------------------------------------------------------------
# Table 2 (rows 45-52, "Kati 1"/floor-1 vertices): update surveyed values and
# drop the last two rows of the table (points "119" and "124"). Before the
# rows are deleted, copy the distinctive "row before total" / "total row"
# formatting down onto the rows that will become those rows once 51 & 52 are
# removed, so the shift lands on the exact same formats as the real edit.
# ---------------------------------------------------------------------------

$ws.Range("B45").Value = 12
$ws.Range("C45").Value = 7510650.6589000002
$ws.Range("D45").Value = 4693484.4923
$ws.Range("E45").Value = 645.94299999999998
$ws.Range("I45").Value = 50.09

$ws.Range("B46").Value = 13
$ws.Range("C46").Value = 7510657.1150000002
$ws.Range("D46").Value = 4693488.3409000002

$ws.Range("B47").Value = 14
$ws.Range("C47").Value = 7510657.0719999997
$ws.Range("D47").Value = 4693488.4129999997
$ws.Range("E47").Value = 646.09199999999998

$ws.Range("B48").Value = 17
$ws.Range("C48").Value = 7510659.7607000005
$ws.Range("D48").Value = 4693489.9873000002

$ws.Range("B49").Value = 18
$ws.Range("C49").Value = 7510662.1574999997
$ws.Range("D49").Value = 4693485.8748000003

$ws.Range("B50").Value = 117
$ws.Range("C50").Value = 7510653.0743000004
$ws.Range("D50").Value = 4693480.4605

# Carry the special formatting from the soon-to-be-deleted rows up onto the
# rows that will take their place (formats only, not values/formulas).
$ws.Range("H51:I51").Copy()
$ws.Range("H49:I49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H52:I52").Copy()
$ws.Range("H50:I50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove rows 51 ("119") and 52 ("124", the old total row); everything below
# shifts up by two rows.
$ws.Rows("51:52").Delete()

# Re-point the merged ranges for the (now shorter) "H"/"I" column of table 2
# and re-establish the "Gjithsej:" total row content/formula on what is now
# row 50.
$ws.Range("H45:H50").UnMerge()
$ws.Range("H45:H48").Merge()
$ws.Range("I45:I50").UnMerge()
$ws.Range("I45:I48").Merge()

$ws.Range("H50").Value = "Gjithsej:"
$ws.Range("I50").Formula = "=SUM(I45:I48)"

# ---------------------------------------------------------------------------
# Sheet view / selection bookkeeping.
# ---------------------------------------------------------------------------

$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("B30:I41").Select()

# ---------------------------------------------------------------------------
# Workbook-level print area.
# ---------------------------------------------------------------------------

$ws.PageSetup.PrintArea = "A1:K56"
